# Fruta / hortaliza, semanal
# Rotates the per-row market data (Fecha, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Unidad de comercializacion,
# Origen, Precio $/Kg, Kg / unidad) across rows 2-13 of Sheet1, matching
# a newer weekly snapshot while leaving the descriptive columns
# (Mercado ID, Mercado, Region, Codreg, Tipo, Producto ID, Producto,
# Categoria ID, Categoria, Variedad) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current ("before") values for column D (Fecha) and the
# contiguous L:T block for every data row before any writes happen, so the
# permutation below reads only original values.
$snapD = @{}
$snapLT = @{}
for ($r = 2; $r -le 13; $r++) {
    $snapD[$r] = $ws.Cells.Item($r, 4).Value2()
    $snapLT[$r] = $ws.Range("L$r`:T$r").Value()
}

# Destination row -> source row (the row whose original D/L:T values move
# into the destination row).
$mapping = @{
    2  = 12
    3  = 4
    4  = 9
    5  = 10
    6  = 11
    7  = 5
    8  = 6
    9  = 13
    10 = 7
    11 = 8
    12 = 3
    13 = 2
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $ws.Range("D$dest").Value2 = $snapD[$src]
    $ws.Range("L$dest`:T$dest").Value = $snapLT[$src]
}
